$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E ("details"/"kccode" shift right by one).
$ws.Columns("E:E").Insert()

# New column E holds the "wcde_code" dimension values.
$ws.Range("E1").Value = "wcde_code"
$ws.Range("E2").Value = 2
$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 3
$ws.Range("E5").Value = 22
$ws.Range("E6").Value = 23
$ws.Range("E7").Value = 4
$ws.Range("E8").Value = 5

# Row 74 previously carried a stray custom row format (fill-only style) -
# drop it so the row matches the formatting of all the other data rows.
$ws.Rows("74:74").ClearFormats()

# Restore the normal selection/view (no scrolled top-left cell, cursor on E9).
$ws.Range("E9").Select()
